$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Cualquiera"

# "1241" must be stored as text (not a number), matching the source diff's
# t="inlineStr" cell type for C4. Temporarily force a text number format so
# Excel doesn't auto-coerce the numeric-looking string into a number, then
# restore the cell's original style so no stray formatting is introduced.
$origStyleC4 = $ws.Range("C4").Style
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1241"
$ws.Range("C4").Style = $origStyleC4

$ws.Range("D4").Value = 16000
$ws.Range("E4").Value = 16
